$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits (header volume/issue number, week date range) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "25"

$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "6/19/2023"
$c9.Characters(47, 9).Text = "6/25/2023"

# --- Type-changing cells: copy format+value from a donor cell holding the exact target content/style ---
$ws.Range("N14").Copy($ws.Range("M14"))
$ws.Range("G14").Copy($ws.Range("G15"))
$ws.Range("H14").Copy($ws.Range("H15"))
$ws.Range("F14").Copy($ws.Range("F16"))
$ws.Range("D14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("D17").Copy($ws.Range("D27"))
$ws.Range("N14").Copy($ws.Range("E27"))

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("L15").Value = -71.428571428571
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = -100
$ws.Range("L16").Value = 50
$ws.Range("M16").Value = -31.25
$ws.Range("N16").Value = -90.09009009009
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 12.5
$ws.Range("I17").Value = 51
$ws.Range("J17").Value = 46
$ws.Range("K17").Value = 10.869565217391
$ws.Range("L17").Value = 75.862068965517
$ws.Range("M17").Value = 82.142857142857
$ws.Range("N17").Value = -13.559322033898
$ws.Range("C18").Value = 5
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 60
$ws.Range("J18").Value = 53
$ws.Range("K18").Value = 13.207547169811
$ws.Range("L18").Value = 62.162162162162
$ws.Range("M18").Value = -6.25
$ws.Range("N18").Value = -90.506329113924
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 45.16129032258
$ws.Range("I19").Value = 228
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = -8.8
$ws.Range("L19").Value = 42.5
$ws.Range("M19").Value = 37.349397590361
$ws.Range("N19").Value = -52
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 71.428571428571
$ws.Range("I20").Value = 60
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = 46.341463414634
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 20
$ws.Range("N20").Value = -96.411483253588
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 69.230769230769
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 56
$ws.Range("H21").Value = 28.571428571428
$ws.Range("I21").Value = 434
$ws.Range("J21").Value = 432
$ws.Range("K21").Value = 0.462962962962
$ws.Range("L21").Value = 52.280701754386
$ws.Range("M21").Value = 21.229050279329
$ws.Range("N21").Value = -86.330708661417
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 13
$ws.Range("K22").Value = -13.333333333333
$ws.Range("L22").Value = 62.5
$ws.Range("M22").Value = 62.5
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -26.666666666666
$ws.Range("F24").Value = 109
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = -4.38596491228
$ws.Range("I24").Value = 766
$ws.Range("J24").Value = 855
$ws.Range("K24").Value = -10.409356725146
$ws.Range("L24").Value = 27.454242928452
$ws.Range("M24").Value = 69.094922737306
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 18
$ws.Range("H25").Value = 80
$ws.Range("I25").Value = 116
$ws.Range("J25").Value = 87
$ws.Range("K25").Value = 33.333333333333
$ws.Range("L25").Value = 36.470588235294
$ws.Range("M25").Value = 9.43396226415
$ws.Range("G26").Value = 1
$ws.Range("L26").Value = -50
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = -26.315789473684
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0
